# The "Förändrad" (Changed) column C for every data row (rows 2-381)
# is bumped from 2023-09-20 (serial 45189) to 2023-09-21 (serial 45190).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

$rng = $ws.Range($ws.Cells.Item(2, 3), $ws.Cells.Item($lastRow, 3))
$rng.Value2 = 45190
